$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 39 / Row 40: Filecoin and Stacks swap places, with refreshed values ---
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("E39").Value = "  +7.16%  "

$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("E40").Value = "  +8.43%  "

# --- Column D ("Price") updates ---
# These look like plain numbers to Excel's normal type-inference
# (e.g. "10.11"), so they would silently be stored as numeric values
# instead of text if assigned directly. Temporarily format the range
# as Text first, so every assignment below is kept as a literal string,
# then restore the default "Normal" style once all values are written
# so no cell is left with a lingering custom number format.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "58.715.80"
$ws.Range("D3").Value = "2.458.72"
$ws.Range("D5").Value = "158.75"
$ws.Range("D6").Value = "496.51"
$ws.Range("D7").Value = "0.615"
$ws.Range("D9").Value = "2.490.14"
$ws.Range("D12").Value = "0.337"
$ws.Range("D14").Value = "2.875.49"
$ws.Range("D15").Value = "58.598.64"
$ws.Range("D16").Value = "21.78"
$ws.Range("D18").Value = "2.464.67"
$ws.Range("D20").Value = "329.95"
$ws.Range("D21").Value = "10.11"
$ws.Range("D24").Value = "58.89"
$ws.Range("D27").Value = "0.994"
$ws.Range("D28").Value = "2.567.63"
$ws.Range("D29").Value = "7.43"
$ws.Range("D30").Value = "0.0₃0808"
$ws.Range("D31").Value = "0.997"
$ws.Range("D33").Value = "151.98"
$ws.Range("D35").Value = "5.47"
$ws.Range("D37").Value = "1.18"
$ws.Range("D39").Value = "1.43"
$ws.Range("D40").Value = "3.66"
$ws.Range("D41").Value = "34.45"
$ws.Range("D42").Value = "285.37"
$ws.Range("D45").Value = "0.992"
$ws.Range("D46").Value = "0.0547"
$ws.Range("D49").Value = "10.27"
$ws.Range("D51").Value = "18.16"

$ws.Range("D2:D51").Style = "Normal"

# --- Column E ("Volume(1h)") updates ---
$ws.Range("E2").Value = "  +5.10%  "
$ws.Range("E3").Value = "  +2.84%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("E5").Value = "  +7.88%  "
$ws.Range("E6").Value = "  +3.67%  "
$ws.Range("E7").Value = "  +23.34%  "
$ws.Range("E8").Value = "  -0.57%  "
$ws.Range("E9").Value = "  +4.16%  "
$ws.Range("E10").Value = "  +15.85%  "
$ws.Range("E11").Value = "  +5.32%  "
$ws.Range("E12").Value = "  +4.19%  "
$ws.Range("E13").Value = "  +1.43%  "
$ws.Range("E14").Value = "  +2.44%  "
$ws.Range("E15").Value = "  +4.11%  "
$ws.Range("E16").Value = "  +7.17%  "
$ws.Range("E17").Value = "  +2.77%  "
$ws.Range("E18").Value = "  +2.98%  "
$ws.Range("E19").Value = "  +5.67%  "
$ws.Range("E20").Value = "  +4.73%  "
$ws.Range("E21").Value = "  +3.60%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("E23").Value = "  +5.17%  "
$ws.Range("E24").Value = "  +3.68%  "
$ws.Range("E25").Value = "  +4.04%  "
$ws.Range("E26").Value = "  +5.21%  "
$ws.Range("E27").Value = "  -0.71%  "
$ws.Range("E28").Value = "  +2.82%  "
$ws.Range("E29").Value = "  +1.91%  "
$ws.Range("E30").Value = "  +4.72%  "
$ws.Range("E31").Value = "  -0.27%  "
$ws.Range("E32").Value = "  +5.92%  "
$ws.Range("E33").Value = "  +2.50%  "
$ws.Range("E34").Value = "  +4.08%  "
$ws.Range("E35").Value = "  +9.41%  "
$ws.Range("E36").Value = "  +8.12%  "
$ws.Range("E37").Value = "  +6.01%  "
$ws.Range("E38").Value = "  +0.78%  "
$ws.Range("E41").Value = "  +3.08%  "
$ws.Range("E42").Value = "  +12.57%  "
$ws.Range("E43").Value = "  +6.60%  "
$ws.Range("E44").Value = "  +4.42%  "
$ws.Range("E45").Value = "  -0.66%  "
$ws.Range("E46").Value = "  +0.78%  "
$ws.Range("E47").Value = "  +5.11%  "
$ws.Range("E48").Value = "  +3.66%  "
$ws.Range("E49").Value = "  +0.50%  "
$ws.Range("E51").Value = "  +6.55%  "
